# WebServiceList.xlsx update — "update the webserice list"
#
# Adds a new row (#10) to the IPCS webservices list describing the
# "recommendation request" service, resizes that row to fit the wrapped
# text, moves the sheet's active-cell selection down to the new row, and
# clears the stale/custom paper-size page-setup hint left over from the
# previous edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 10 data -------------------------------------------------
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "recommendation request,send email to school to recommend the service"

# Row grows tall enough to show the wrapped remark text.
$ws.Rows.Item(12).RowHeight = 60

# --- Selection moves to the newly edited cell -------------------------
[void]$ws.Range("B13").Select()

# --- Drop the stale page setup (paper size / dpi hints) ----------------
$ws.PageSetup.PaperSize = 0
